$d = $word.ActiveDocument

$d.Content.Find.Execute("391×5=1955", $true, $true, $false, $false, $false, $true, 1, $false, "844×6=5064", 2) | Out-Null
$d.Content.Find.Execute("329×3=987", $true, $true, $false, $false, $false, $true, 1, $false, "951×8=7608", 2) | Out-Null
$d.Content.Find.Execute("767×3=2301", $true, $true, $false, $false, $false, $true, 1, $false, "697×6=4182", 2) | Out-Null
$d.Content.Find.Execute("267×2=534", $true, $true, $false, $false, $false, $true, 1, $false, "521×8=4168", 2) | Out-Null
$d.Content.Find.Execute("281×7=1967", $true, $true, $false, $false, $false, $true, 1, $false, "168×6=1008", 2) | Out-Null
$d.Content.Find.Execute("144×8=1152", $true, $true, $false, $false, $false, $true, 1, $false, "497×9=4473", 2) | Out-Null
$d.Content.Find.Execute("463×8=3704", $true, $true, $false, $false, $false, $true, 1, $false, "116×3=348", 2) | Out-Null
$d.Content.Find.Execute("321×9=2889", $true, $true, $false, $false, $false, $true, 1, $false, "867×6=5202", 2) | Out-Null
$d.Content.Find.Execute("595×3=1785", $true, $true, $false, $false, $false, $true, 1, $false, "136×8=1088", 2) | Out-Null
$d.Content.Find.Execute("363×2=726", $true, $true, $false, $false, $false, $true, 1, $false, "191×6=1146", 2) | Out-Null
$d.Content.Find.Execute("629×6=3774", $true, $true, $false, $false, $false, $true, 1, $false, "511×3=1533", 2) | Out-Null
$d.Content.Find.Execute("139×7=973", $true, $true, $false, $false, $false, $true, 1, $false, "240×4=960", 2) | Out-Null
$d.Content.Find.Execute("682×3=2046", $true, $true, $false, $false, $false, $true, 1, $false, "818×5=4090", 2) | Out-Null
$d.Content.Find.Execute("719×2=1438", $true, $true, $false, $false, $false, $true, 1, $false, "381×4=1524", 2) | Out-Null
$d.Content.Find.Execute("263×7=1841", $true, $true, $false, $false, $false, $true, 1, $false, "507×5=2535", 2) | Out-Null
$d.Content.Find.Execute("301×4=1204", $true, $true, $false, $false, $false, $true, 1, $false, "474×4=1896", 2) | Out-Null
$d.Content.Find.Execute("139×6=834", $true, $true, $false, $false, $false, $true, 1, $false, "463×7=3241", 2) | Out-Null
$d.Content.Find.Execute("383×3=1149", $true, $true, $false, $false, $false, $true, 1, $false, "757×9=6813", 2) | Out-Null
$d.Content.Find.Execute("734×2=1468", $true, $true, $false, $false, $false, $true, 1, $false, "528×3=1584", 2) | Out-Null
$d.Content.Find.Execute("733×7=5131", $true, $true, $false, $false, $false, $true, 1, $false, "444×3=1332", 2) | Out-Null
$d.Content.Find.Execute("669×4=2676", $true, $true, $false, $false, $false, $true, 1, $false, "619×3=1857", 2) | Out-Null
$d.Content.Find.Execute("635×9=5715", $true, $true, $false, $false, $false, $true, 1, $false, "363×9=3267", 2) | Out-Null
$d.Content.Find.Execute("423×5=2115", $true, $true, $false, $false, $false, $true, 1, $false, "685×3=2055", 2) | Out-Null
$d.Content.Find.Execute("387×4=1548", $true, $true, $false, $false, $false, $true, 1, $false, "522×5=2610", 2) | Out-Null
$d.Content.Find.Execute("849×3=2547", $true, $true, $false, $false, $false, $true, 1, $false, "860×9=7740", 2) | Out-Null
